$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edit: update the seed "current value" (A4) and "dollars added per
#     purchase" (B4) inputs. Row-2 KPI formulas (I2/J2/K2/L2) recalc off these. ---
$ws.Range("A4").Value = 2.053
$ws.Range("B4").Value = 10

# --- Fill colour touch-up -------------------------------------------------
# D4 previously shared the bright-green swatch (RGB D4EA6B / BBE33D pair);
# recolour it to the lighter green pairing (D4EA6B / E8F2A1).
$ws.Range("D4").Interior.Color = 7072468
$ws.Range("D4").Interior.PatternColor = 10613480

# A4:C4 drop the leftover green fill entirely (back to no fill), matching
# the already-unfilled cells below them (A5:C8).
$ws.Range("A4:C4").Interior.ColorIndex = -4142

# --- Used range / selection bookkeeping -----------------------------------
# Touch row 11 so the sheet's used range (and <dimension>) extends to it,
# then clear it back out so it stays an empty trailing row.
$ws.Cells.Item(11, 1).Value = "x"
$ws.Cells.Item(11, 1).ClearContents()

# Reflect the new selection left in the sheet view.
$ws.Range("A5:D20").Select()
